$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above row 7 (shifts existing rows 7.. down by one).
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the added BOM line (connector terminals).
$ws.Cells.Item(7, 2).Value = 180
$ws.Cells.Item(7, 3).Value = "S9473CT-ND"
$ws.Cells.Item(7, 6).Value = "Connector Terminals"
$ws.Cells.Item(7, 8).Value = 53

# Re-apply the style used by ordinary data rows to the new row's cells.
$ws.Cells.Item(7, 2).Style = $ws.Cells.Item(8, 2).Style
$ws.Cells.Item(7, 3).Style = $ws.Cells.Item(8, 3).Style
$ws.Cells.Item(7, 6).Style = $ws.Cells.Item(8, 6).Style
$ws.Cells.Item(7, 8).Style = $ws.Cells.Item(8, 8).Style

# Move the view/selection to match the edited state (scrolled down to the
# newly-added last data row, selection on it).
$ws.Application.GoTo($ws.Range("B72"))
$excel.ActiveWindow.ScrollRow = 49
$ws.Range("B72").Select()
